# Insert a new data row into the weekly Hortaliza/Zanahoria log.
# This pushes the existing rows 281..299 down to 282..300 (preserving all
# their values/format) and fills the freshly inserted row 281 with the new
# weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("281:281").Insert()

$ws.Cells.Item(281, 1).Value = 3
$ws.Cells.Item(281, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(281, 3).Value = "Coquimbo"
$ws.Cells.Item(281, 4).Value = 44610
$ws.Cells.Item(281, 5).Value = 5
$ws.Cells.Item(281, 6).Value = 100114013
$ws.Cells.Item(281, 7).Value = "Zanahoria"
$ws.Cells.Item(281, 8).Value = "Sin especificar"
$ws.Cells.Item(281, 9).Value = "Primera"
$ws.Cells.Item(281, 10).Value = 230
$ws.Cells.Item(281, 11).Value = 9500
$ws.Cells.Item(281, 12).Value = 10000
$ws.Cells.Item(281, 13).Value = 9739
$ws.Cells.Item(281, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(281, 15).Value = "Chillán"
$ws.Cells.Item(281, 16).Value = 487
$ws.Cells.Item(281, 17).Value = 20
$ws.Cells.Item(281, 18).Value = "Hortaliza"
